$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 45000

# Row 3
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 45000

# Row 4
$ws.Range("B4").Value = 12500
$ws.Range("C4").Value = 10000
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 22500

# Row 5
$ws.Range("B5").Value = 12500
$ws.Range("C5").Value = 10000
$ws.Range("E5").Value = 22500

# Row 6
$ws.Range("B6").Value = 12500
$ws.Range("E6").Value = 32500

# Row 10
$ws.Range("B10").Value = 25000
$ws.Range("C10").Value = 20000
$ws.Range("D10").Value = 20000
$ws.Range("E10").Value = 65000

# Row 11
$ws.Range("B11").Value = 25000
$ws.Range("C11").Value = 20000
$ws.Range("D11").Value = 20000
$ws.Range("E11").Value = 65000

# Row 12
$ws.Range("D12").Value = 20000
$ws.Range("E12").Value = 65000

# Row 13
$ws.Range("D13").Value = 20000
$ws.Range("E13").Value = 65000

# Row 14
$ws.Range("D14").Value = 20000
$ws.Range("E14").Value = 65000
